$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1183.9803
$ws.Range("J17").Value = 1183.9803
$ws.Range("L17").Value = 3551.9409
$ws.Range("N17").Value = -3887.9409

$ws.Range("H98").Value = 2697.282
$ws.Range("I98").Value = 2863.1765
$ws.Range("K98").Value = 2863.1765
$ws.Range("M98").Value = -1365.1765

$ws.Range("H101").Value = 819.6
$ws.Range("I101").Value = 699.3333
$ws.Range("J101").Value = 1000
$ws.Range("K101").Value = 2097.9999
$ws.Range("L101").Value = 3000
$ws.Range("M101").Value = -475.9998999999998
$ws.Range("N101").Value = -6244

$ws.Range("H112").Value = 1967.75
$ws.Range("I112").Value = 990.9091
$ws.Range("J112").Value = 2258.162
$ws.Range("K112").Value = 2972.7273
$ws.Range("L112").Value = 6774.485999999999
$ws.Range("M112").Value = -1864.7273
$ws.Range("N112").Value = -8990.485999999999

$ws.Range("H122").Value = 2697.282
$ws.Range("I122").Value = 2863.1765
$ws.Range("K122").Value = 8589.529500000001
$ws.Range("M122").Value = -6139.529500000001

$ws.Range("H129").Value = 538.4
$ws.Range("J129").Value = 887.9091
$ws.Range("L129").Value = 2663.7273
$ws.Range("N129").Value = -12663.7273

$ws.Range("H132").Value = 5054872.5
$ws.Range("I132").Value = 6948344
$ws.Range("J132").Value = 5613.8335
$ws.Range("K132").Value = 20845032
$ws.Range("L132").Value = 16841.5005
$ws.Range("M132").Value = -20842502
$ws.Range("N132").Value = -21901.5005

$ws.Range("H138").Value = 1289.33
$ws.Range("I138").Value = 685.2857
$ws.Range("J138").Value = 1614.5846
$ws.Range("K138").Value = 2055.8571
$ws.Range("L138").Value = 4843.7538
$ws.Range("M138").Value = 3084.1429
$ws.Range("N138").Value = -15123.7538

$ws.Range("H141").Value = 659.7059
$ws.Range("I141").Value = 575.9375
$ws.Range("K141").Value = 1727.8125
$ws.Range("M141").Value = 3452.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 32259312
$ws.Range("I61").Value = 43479160
$ws.Range("J61").Value = 2250
$ws.Range("K61").Value = 43479160
$ws.Range("L61").Value = 2250
$ws.Range("M61").Value = -43478948
$ws.Range("N61").Value = -2674

$ws.Range("H110").Value = 1364.9615
$ws.Range("I110").Value = 914.1539
$ws.Range("K110").Value = 914.1539
$ws.Range("M110").Value = 1130.8461

$ws.Range("H132").Value = 1629.0952
$ws.Range("I132").Value = 1345.9678
$ws.Range("J132").Value = 2427
$ws.Range("K132").Value = 4037.9034
$ws.Range("L132").Value = 7281
$ws.Range("M132").Value = -1507.9034
$ws.Range("N132").Value = -12341

$ws.Range("H136").Value = 32259312
$ws.Range("I136").Value = 43479160
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 130437480
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -130434930
$ws.Range("N136").Value = -11850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1013
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1013
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = ""
$ws.Range("M22").Value = 1013
$ws.Range("N22").Value = -1359

$ws.Range("H99").Value = 58824696
$ws.Range("I99").Value = 83334360
$ws.Range("K99").Value = 83334360
$ws.Range("M99").Value = -83332862

$ws.Range("H134").Value = 3014.5095
$ws.Range("I134").Value = 933.9535
$ws.Range("K134").Value = 2801.8605
$ws.Range("M134").Value = -266.8604999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1882.5769
$ws.Range("I31").Value = 1882.5769
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1882.5769
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -1587.5769

$ws.Range("H34").Value = 1882.5769
$ws.Range("I34").Value = 1882.5769
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1882.5769
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -1680.5769

$ws.Range("H58").Value = 1020.6667
$ws.Range("I58").Value = 970.27026
$ws.Range("K58").Value = 970.27026
$ws.Range("M58").Value = -767.27026

$ws.Range("H132").Value = 1516.3191
$ws.Range("I132").Value = 1154.75
$ws.Range("K132").Value = 3464.25
$ws.Range("M132").Value = -934.25

$ws.Range("H134").Value = 14286792
$ws.Range("I134").Value = 988.2083
$ws.Range("J134").Value = 45455820
$ws.Range("K134").Value = 2964.6249
$ws.Range("L134").Value = 136367460
$ws.Range("M134").Value = -429.6248999999998
$ws.Range("N134").Value = -136372530

$ws.Range("H136").Value = 1020.6667
$ws.Range("I136").Value = 970.27026
$ws.Range("K136").Value = 2910.81078
$ws.Range("M136").Value = -360.8107799999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1130.8572
$ws.Range("I98").Value = 1469.2222
$ws.Range("J98").Value = 521.8
$ws.Range("K98").Value = 4407.6666
$ws.Range("L98").Value = 1565.4
$ws.Range("M98").Value = -2909.6666
$ws.Range("N98").Value = -4561.4

$ws.Range("H131").Value = 27030746
$ws.Range("I131").Value = 125000500
$ws.Range("J131").Value = 4606.6206
$ws.Range("K131").Value = 375001500
$ws.Range("L131").Value = 13819.8618
$ws.Range("M131").Value = -374996460
$ws.Range("N131").Value = -23899.8618

$ws.Range("H132").Value = 1668.75
$ws.Range("I132").Value = 1154.5454
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 10390.9086
$ws.Range("L132").Value = 25200
$ws.Range("M132").Value = -7860.908599999999
$ws.Range("N132").Value = -30260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 30148.834
$ws.Range("J86").Value = 30148.834
$ws.Range("L86").Value = 30148.834
$ws.Range("N86").Value = -32520.834

$ws.Range("H89").Value = 30148.834
$ws.Range("J89").Value = 30148.834
$ws.Range("L89").Value = 90446.50199999999
$ws.Range("N89").Value = -102302.502

$ws.Range("H132").Value = 2034.1666
$ws.Range("I132").Value = 1586.6333
$ws.Range("J132").Value = 3153
$ws.Range("K132").Value = 4759.8999
$ws.Range("L132").Value = 9459
$ws.Range("M132").Value = -2229.8999
$ws.Range("N132").Value = -14519

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8600
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = ""

$ws.Range("H93").Value = 606.1111
$ws.Range("I93").Value = 550.8570999999999
$ws.Range("J93").Value = 799.5
$ws.Range("K93").Value = 550.8570999999999
$ws.Range("L93").Value = 799.5
$ws.Range("M93").Value = 697.1429000000001
$ws.Range("N93").Value = -3295.5

$ws.Range("H132").Value = 25658.357
$ws.Range("I132").Value = 1152.2727
$ws.Range("K132").Value = 3456.8181
$ws.Range("M132").Value = -926.8181

$ws.Range("H136").Value = 2327.2727
$ws.Range("I136").Value = 2288.889
$ws.Range("K136").Value = 6866.667
$ws.Range("M136").Value = -4316.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 550.75
$ws.Range("J81").Value = 701
$ws.Range("L81").Value = 1402
$ws.Range("N81").Value = -3524

$ws.Range("H84").Value = 550.75
$ws.Range("J84").Value = 701
$ws.Range("L84").Value = 7010
$ws.Range("N84").Value = -17618

$ws.Range("H132").Value = 1689.5652
$ws.Range("I132").Value = 1182.5714
$ws.Range("J132").Value = 2478.2222
$ws.Range("K132").Value = 3547.7142
$ws.Range("L132").Value = 7434.6666
$ws.Range("M132").Value = -1017.7142
$ws.Range("N132").Value = -12494.6666

$ws.Range("H136").Value = 599.1429000000001
$ws.Range("I136").Value = 518.2759
$ws.Range("K136").Value = 1554.8277
$ws.Range("M136").Value = 995.1723000000002
